# ============================================================================
# Edit script: apply the "08-03-2024" batch of new transactions across the
# HỆ THỐNG 8 - 2024 workbook, and rebuild the "LỢI NHUẬN" summary template.
# ============================================================================

$wb = $excel.ActiveWorkbook

function Set-Text($ws, $row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
}

function Set-Num($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Clear-Cell($ws, $row, $col) {
    $ws.Cells.Item($row, $col).Value = ""
}

# ----------------------------------------------------------------------------
# Sheet 1: CHI TIẾT DOANH THU
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 6 was a half-empty data row (SÓC TRĂNG) - fill in the missing details
Set-Text $ws1 6 5  "Cắt mí"
Set-Text $ws1 6 6  "đường thị út"
Set-Text $ws1 6 7  "Cá nhân"
Set-Text $ws1 6 8  "Lê Đình Hậu"
Set-Num  $ws1 6 9  6000000
Set-Num  $ws1 6 12 6000000
Set-Text $ws1 6 13 "Nguyễn Hoàng Yến Quyên"
Set-Num  $ws1 6 15 6000000
Set-Num  $ws1 6 17 6000000
Set-Text $ws1 6 19 "Kha Như Huỳnh "
Set-Num  $ws1 6 21 50000
Set-Num  $ws1 6 22 0

# Row 8 used to hold the "Tổng" summary - it becomes a brand-new data row
Set-Text $ws1 8 1  "08-03-2024"
Set-Text $ws1 8 2  "HD-LUXURY"
Set-Num  $ws1 8 3  620
Set-Text $ws1 8 4  "CẦN THƠ"
Set-Text $ws1 8 5  "Nâng mũi"
Set-Text $ws1 8 6  "Trần Thị Thanh Nhàn"
Set-Text $ws1 8 7  "Cá nhân"
Set-Text $ws1 8 8  "Phạm Thanh Hoàng"
Set-Num  $ws1 8 9  28000000
Clear-Cell $ws1 8 10
Clear-Cell $ws1 8 11
Set-Num  $ws1 8 12 28000000
Set-Text $ws1 8 13 "Phạm Thanh Hoàng"
Clear-Cell $ws1 8 14
Set-Num  $ws1 8 15 28000000
Set-Num  $ws1 8 16 0
Set-Num  $ws1 8 17 28000000
Set-Num  $ws1 8 18 0
Set-Text $ws1 8 19 "Lâm Hoàng Phú"
Clear-Cell $ws1 8 20
Set-Num  $ws1 8 21 100000
Set-Num  $ws1 8 22 50000

# Row 9 is an entirely new data row
Set-Text $ws1 9 1  "08-03-2024"
Set-Text $ws1 9 2  "HD-LUXURY"
Set-Num  $ws1 9 3  621
Set-Text $ws1 9 4  "CẦN THƠ"
Set-Text $ws1 9 5  "Nâng mũi"
Set-Text $ws1 9 6  "Trần Thị Ngọc Dung"
Set-Text $ws1 9 7  "Cá nhân"
Set-Text $ws1 9 8  "Lâm Thị Mỹ Hằng"
Set-Num  $ws1 9 9  9000000
Clear-Cell $ws1 9 10
Clear-Cell $ws1 9 11
Set-Num  $ws1 9 12 9000000
Set-Text $ws1 9 13 "Phạm Thanh Hoàng"
Clear-Cell $ws1 9 14
Set-Num  $ws1 9 15 9000000
Set-Num  $ws1 9 16 0
Set-Num  $ws1 9 17 9000000
Set-Num  $ws1 9 18 0
Set-Text $ws1 9 19 "Lâm Hoàng Phú"
Clear-Cell $ws1 9 20
Set-Num  $ws1 9 21 100000
Set-Num  $ws1 9 22 50000

# Row 10 is the new "Tổng" row (shifted down from row 8)
Clear-Cell $ws1 10 1
Set-Text $ws1 10 2  "Tổng"
Set-Num  $ws1 10 3  8
Clear-Cell $ws1 10 4
Clear-Cell $ws1 10 5
Clear-Cell $ws1 10 6
Clear-Cell $ws1 10 7
Clear-Cell $ws1 10 8
Set-Num  $ws1 10 9  65100000
Clear-Cell $ws1 10 10
Set-Num  $ws1 10 11 6000000
Set-Num  $ws1 10 12 71100000
Clear-Cell $ws1 10 13
Clear-Cell $ws1 10 14
Set-Num  $ws1 10 15 68100000
Set-Num  $ws1 10 16 0
Set-Num  $ws1 10 17 68100000
Set-Num  $ws1 10 18 3000000
Clear-Cell $ws1 10 19
Clear-Cell $ws1 10 20
Set-Num  $ws1 10 21 550000
Set-Num  $ws1 10 22 100000

Write-Host "Sheet1 done"

# ----------------------------------------------------------------------------
# Sheet 2: CHI TIẾT VỀ THU NỢ
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New column inserted before "Lượng thu": "Ngày thực hiện"
Set-Text $ws2 1 6 "Ngày thực hiện"
Set-Text $ws2 1 7 "Lượng thu"

# Row 2 used to be "Tổng" - now becomes the first debt-collection record
Set-Text $ws2 2 1 "TN"
Set-Num  $ws2 2 2 171
Set-Text $ws2 2 3 "08-03-2024"
Set-Text $ws2 2 4 "CẦN THƠ"
Set-Text $ws2 2 5 "HD-LUXURY-190"
Set-Text $ws2 2 6 "2024-04-16"
Set-Num  $ws2 2 7 1000000

# Row 3 is a new record
Set-Text $ws2 3 1 "TN"
Set-Num  $ws2 3 2 172
Set-Text $ws2 3 3 "08-03-2024"
Set-Text $ws2 3 4 "CẦN THƠ"
Set-Text $ws2 3 5 "HD-LUXURY-191"
Set-Text $ws2 3 6 "2024-04-16"
Set-Num  $ws2 3 7 4000000

# Row 4 is a new record
Set-Text $ws2 4 1 "TN"
Set-Num  $ws2 4 2 173
Set-Text $ws2 4 3 "08-03-2024"
Set-Text $ws2 4 4 "CẦN THƠ"
Set-Text $ws2 4 5 "HD-LUXURY-587"
Set-Text $ws2 4 6 "2024-07-24"
Set-Num  $ws2 4 7 10000000

# Row 5 is the new "Tổng" row (shifted down from row 2)
Set-Text $ws2 5 1 "Tổng"
Set-Num  $ws2 5 2 3
Clear-Cell $ws2 5 3
Clear-Cell $ws2 5 4
Clear-Cell $ws2 5 5
Clear-Cell $ws2 5 6
Set-Num  $ws2 5 7 15000000

Write-Host "Sheet2 done"

# ----------------------------------------------------------------------------
# Sheet 3: CHI TIẾT CHI TIÊU
# ----------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# 5 new expense rows inserted before the old "Tổng" row (which moves 11 -> 16)
Set-Text $ws3 11 1 "CT"
Set-Num  $ws3 11 2 750
Set-Text $ws3 11 3 "08-03-2024"
Set-Text $ws3 11 4 "CẦN THƠ"
Set-Text $ws3 11 5 "Ứng Lương"
Set-Num  $ws3 11 6 1000000

Set-Text $ws3 12 1 "CT"
Set-Num  $ws3 12 2 751
Set-Text $ws3 12 3 "08-03-2024"
Set-Text $ws3 12 4 "CẦN THƠ"
Set-Text $ws3 12 5 "Chi Phí Sinh Hoạt Tại Cơ Sở"
Set-Num  $ws3 12 6 2230000

Set-Text $ws3 13 1 "CT"
Set-Num  $ws3 13 2 752
Set-Text $ws3 13 3 "08-03-2024"
Set-Text $ws3 13 4 "SÓC TRĂNG"
Set-Text $ws3 13 5 "Chi Phí Sinh Hoạt Tại Cơ Sở"
Set-Num  $ws3 13 6 180000

Set-Text $ws3 14 1 "CT"
Set-Num  $ws3 14 2 753
Set-Text $ws3 14 3 "08-03-2024"
Set-Text $ws3 14 4 "SÓC TRĂNG"
Set-Text $ws3 14 5 "Chi Phí Vận Hành"
Set-Num  $ws3 14 6 500000

Set-Text $ws3 15 1 "CT"
Set-Num  $ws3 15 2 754
Set-Text $ws3 15 3 "08-03-2024"
Set-Text $ws3 15 4 "SÓC TRĂNG"
Set-Text $ws3 15 5 "Trang thiết bị Y Tế"
Set-Num  $ws3 15 6 700000

# Row 16 is the new "Tổng" row (shifted down from row 11)
Set-Text $ws3 16 1 "Tổng"
Set-Num  $ws3 16 2 14
Clear-Cell $ws3 16 3
Clear-Cell $ws3 16 4
Clear-Cell $ws3 16 5
Set-Num  $ws3 16 6 14010000

Write-Host "Sheet3 done"

# ----------------------------------------------------------------------------
# Sheet 4: DOANH SỐ CÁ NHÂN
# ----------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-Num $ws4 3 6  2        # Kha Như Huỳnh - Số lần phụ phẫu 1
Set-Num $ws4 3 7  150000   # Kha Như Huỳnh - Công phụ phẫu 1

Set-Num $ws4 4 6  3        # Lâm Hoàng Phú - Số lần phụ phẫu 1
Set-Num $ws4 4 7  250000   # Lâm Hoàng Phú - Công phụ phẫu 1

Set-Num $ws4 5 2  9000000   # Lâm Thị Mỹ Hằng - Tổng đơn giá sale vòng 1
Set-Num $ws4 5 10 159800000 # Lâm Thị Mỹ Hằng - Doanh số thu nợ

Set-Num $ws4 7 10 36000000  # Lê Văn Linh - Doanh số thu nợ

Set-Num $ws4 8 2  14000000  # Lê Đình Hậu - Tổng đơn giá sale vòng 1

Set-Num $ws4 9 4  25100000  # Nguyễn Hoàng Yến Quyên - Doanh số đơn 1 bác sĩ

Set-Num $ws4 11 2 28000000  # Phạm Thanh Hoàng - Tổng đơn giá sale vòng 1
Set-Num $ws4 11 4 37000000  # Phạm Thanh Hoàng - Doanh số đơn 1 bác sĩ

Set-Num $ws4 17 2 65100000  # Tổng - Tổng đơn giá sale vòng 1
Set-Num $ws4 17 4 68100000  # Tổng - Doanh số đơn 1 bác sĩ
Set-Num $ws4 17 6 7         # Tổng - Số lần phụ phẫu 1
Set-Num $ws4 17 7 500000    # Tổng - Công phụ phẫu 1
Set-Num $ws4 17 10 460186000 # Tổng - Doanh số thu nợ

Write-Host "Sheet4 done"

# ----------------------------------------------------------------------------
# Sheet 5: CHI TIÊU TỔNG HỢP
# ----------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-Num $ws5 3 2 3090000   # Chi Phí Sinh Hoạt Tại Cơ Sở

Set-Num $ws5 4 2 1330000   # Chi Phí Vận Hành

Set-Num $ws5 5 2 3090000   # Trang thiết bị Y Tế

# Row 6 is a brand-new "Ứng Lương" category; "Blank"/"Tổng cộng" rows shift down by 1
Set-Text $ws5 6 1 "Ứng Lương"
Set-Num  $ws5 6 2 1000000

Set-Text $ws5 7 1 "Blank"
Set-Num  $ws5 7 2 5000000

Set-Text $ws5 8 1 "Tổng cộng"
Set-Num  $ws5 8 2 14010000

Write-Host "Sheet5 done"

# ----------------------------------------------------------------------------
# Sheet 6: LŨY KẾ NGÀY
# ----------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-Num $ws6 3 2 26000000  # 08-02-2024 - Đơn giá
Set-Num $ws6 3 3 23000000  # 08-02-2024 - Thanh toán lần đầu
Set-Num $ws6 3 7 19220000  # 08-02-2024 - Lũy kế ngày

# Row 4 is the brand-new 08-03-2024 entry; "Tổng" shifts down 4 -> 5
Set-Text $ws6 4 1 "08-03-2024"
Set-Num  $ws6 4 2 37000000
Set-Num  $ws6 4 3 37000000
Set-Num  $ws6 4 4 2
Set-Num  $ws6 4 5 15000000
Set-Num  $ws6 4 6 4610000
Set-Num  $ws6 4 7 47390000

Set-Text $ws6 5 1 "Tổng"
Set-Num  $ws6 5 2 71100000
Set-Num  $ws6 5 3 68100000
Set-Num  $ws6 5 4 8
Set-Num  $ws6 5 5 15000000
Set-Num  $ws6 5 6 14010000
Set-Num  $ws6 5 7 69090000

Write-Host "Sheet6 done"

# ----------------------------------------------------------------------------
# Sheet 7: QUỸ LƯƠNG
# ----------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Row 4 (NV-29, Lâm Hoàng Phú): CẦN THƠ pay increases
Set-Num $ws7 4 3 482142.8571428572
Set-Num $ws7 4 6 482142.8571428572

# A new row 9 (NV-5, Nguyễn Hoàng Yến Quyên) is inserted; every row from the
# old row 9 onward shifts down by one. Two more rows get inserted further
# down (new NV-22 at 13, new NV-10 at 18), so build the whole row9..row23
# block explicitly in its final, post-insert order.

Set-Text $ws7 9 1 "NV-5"
Set-Text $ws7 9 2 "Nguyễn Hoàng Yến Quyên"
Set-Num  $ws7 9 3 1289428.571428571
Set-Num  $ws7 9 4 1385714.285714286
Set-Num  $ws7 9 5 600000
Set-Num  $ws7 9 6 3275142.857142857

Set-Text $ws7 10 1 "NV-6"
Set-Text $ws7 10 2 "Lâm Thị Mỹ Hằng"
Set-Num  $ws7 10 3 2378095.238095238
Set-Num  $ws7 10 4 238095.2380952381
Set-Num  $ws7 10 5 838095.2380952381
Set-Num  $ws7 10 6 3454285.714285714

Set-Text $ws7 11 1 "NV-7"
Set-Text $ws7 11 2 "Phạm Thanh Hoàng"
Set-Num  $ws7 11 3 7771071.428571429
Set-Num  $ws7 11 4 53571.42857142857
Set-Num  $ws7 11 5 53571.42857142857
Set-Num  $ws7 11 6 7878214.285714285

Set-Text $ws7 12 1 "NV-9"
Set-Text $ws7 12 2 "Lê Văn Linh"
Set-Num  $ws7 12 3 2141428.571428571
Set-Num  $ws7 12 4 714285.7142857143
Set-Num  $ws7 12 5 1071428.571428571
Set-Num  $ws7 12 6 3927142.857142857

Set-Text $ws7 13 1 "NV-22"
Set-Text $ws7 13 2 "Nguyễn Phúc Nam"
Set-Num  $ws7 13 3 0
Set-Num  $ws7 13 4 1957142.857142857
Set-Num  $ws7 13 5 0
Set-Num  $ws7 13 6 1957142.857142857

Set-Text $ws7 14 1 "NV-23"
Set-Text $ws7 14 2 "Lê Hoàng Thanh"
Set-Num  $ws7 14 3 0
Set-Num  $ws7 14 4 535714.2857142857
Set-Num  $ws7 14 5 0
Set-Num  $ws7 14 6 535714.2857142857

Set-Text $ws7 15 1 "NV-30"
Set-Text $ws7 15 2 "Đào Vương Anh"
Set-Num  $ws7 15 3 0
Set-Num  $ws7 15 4 528571.4285714286
Set-Num  $ws7 15 5 0
Set-Num  $ws7 15 6 528571.4285714286

Set-Text $ws7 16 1 "NV-36"
Set-Text $ws7 16 2 "Đặng Ngọc Mai"
Set-Num  $ws7 16 3 0
Set-Num  $ws7 16 4 640714.2857142857
Set-Num  $ws7 16 5 0
Set-Num  $ws7 16 6 640714.2857142857

Set-Text $ws7 17 1 "NV-40"
Set-Text $ws7 17 2 "Sang sang"
Set-Num  $ws7 17 3 0
Set-Num  $ws7 17 4 355357.1428571428
Set-Num  $ws7 17 5 0
Set-Num  $ws7 17 6 355357.1428571428

Set-Text $ws7 18 1 "NV-10"
Set-Text $ws7 18 2 "Lê Đình Hậu"
Set-Num  $ws7 18 3 0
Set-Num  $ws7 18 4 0
Set-Num  $ws7 18 5 1380000
Set-Num  $ws7 18 6 1380000

Set-Text $ws7 19 1 "NV-16"
Set-Text $ws7 19 2 "Kha Như Huỳnh "
Set-Num  $ws7 19 3 0
Set-Num  $ws7 19 4 0
Set-Num  $ws7 19 5 150000
Set-Num  $ws7 19 6 150000

Set-Text $ws7 20 1 "NV-26"
Set-Text $ws7 20 2 "Trần Khánh Hiệp"
Set-Num  $ws7 20 3 0
Set-Num  $ws7 20 4 0
Set-Num  $ws7 20 5 0
Set-Num  $ws7 20 6 0

Set-Text $ws7 21 1 "NV-27"
Set-Text $ws7 21 2 "Cô Siêng giúp Việc"
Set-Num  $ws7 21 3 0
Set-Num  $ws7 21 4 0
Set-Num  $ws7 21 5 0
Set-Num  $ws7 21 6 0

Set-Text $ws7 22 1 "NV-38"
Set-Text $ws7 22 2 "Lê Thị Ngọc Mi"
Set-Num  $ws7 22 3 0
Set-Num  $ws7 22 4 0
Set-Num  $ws7 22 5 0
Set-Num  $ws7 22 6 0

Set-Text $ws7 23 1 "Tổng lương"
Clear-Cell $ws7 23 2
Set-Num  $ws7 23 3 16225023.80952381
Set-Num  $ws7 23 4 6766309.523809524
Set-Num  $ws7 23 5 4521666.666666666
Set-Num  $ws7 23 6 27513000

Write-Host "Sheet7 done"
